$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ID column (A2:A6) with the new example IDs
$ws.Range("A2").Value = "ID6"
$ws.Range("A3").Value = "ID7"
$ws.Range("A4").Value = "ID8"
$ws.Range("A5").Value = "ID9"
$ws.Range("A6").Value = "ID10"

# Match the saved selection state
$ws.Range("A6").Select()
